# Trade #3 closed at 2026-02-17 15:13:30 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook after a new (break-even) trade
# closes:
#   - Summary sheet: capital/P&L/trade-count/win-rate roll-ups
#   - Strategy Status sheet: MarketMaking strategy row roll-ups
#   - All Trades / MarketMaking sheets: append the new trade as row 4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200        # Current Capital
$summary.Range("B4").Value = 0           # Total P&L $
$summary.Range("B5").Value = 0           # Total P&L %
$summary.Range("B6").Value = 3           # Total Trades
$summary.Range("B7").Value = 2           # Winning Trades
$summary.Range("B9").Value = 66.67       # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100          # Capital
$status.Range("D4").Value = 3            # Trades
$status.Range("E4").Value = 0            # P&L $
$status.Range("F4").Value = 0            # P&L %
$status.Range("G4").Value = 66.67        # Win Rate %

# ---------------------------------------------------------------------------
# All Trades / MarketMaking - append the new closed trade as row 4.
#
# Row 3 is copied down to row 4 first so the new row inherits the existing
# cell formatting/type (in particular keeping the Date/Time columns as plain
# text instead of Excel re-interpreting "2026-02-17" as a date serial), then
# just the cells that differ for this trade are overwritten.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A3:Q3").Copy($ws.Range("A4:Q4"))

    $ws.Range("A4").Value = 3                                          # Trade #
    # B4 Date, D4 Strategy, E4 Side, H4 Status, L4/M4 slippage, N4
    # Confidence, O4 Entry Reason, P4 Exit Reason are unchanged from row 3.
    $ws.Range("C4").Value = "15:13:24"                                 # Time
    $ws.Range("F4").Value = 0.82                                       # Entry Price
    $ws.Range("G4").Value = 0.86                                       # Exit Price
    $ws.Range("I4").Value = 4.878                                      # P&L %
    $ws.Range("J4").Value = 0.04                                       # P&L $
    $ws.Range("K4").Value = 100                                        # Capital After
    $ws.Range("Q4").Value = 0.14                                       # Duration (min)
}
